$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45180 -> 45181, i.e. 2023-09-11 -> 2023-09-12) for every data row
# (rows 2 through 307).
$ws.Range("C2:C307").Value = 45181
